$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mis-encoded "1KΩ" value (was "1KÎ©") for R2 / R1 resistor rows ---
$ws.Range("B10").Value = "1KΩ"
$ws.Range("B11").Value = "1KΩ"

# --- Swap the Description/Comments headers (E1 <-> F1) ---
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Comments"

# --- Move the existing "Comments" text out to column F ---
$ws.Range("F3").Value  = "Included breakout board"
$ws.Range("F5").Value  = "2 PACK"
$ws.Range("F6").Value  = "included in 2 pack above"
$ws.Range("F7").Value  = "Included breakout board"
$ws.Range("F10").Value = "Must order at least 10 at $0.099 each"
$ws.Range("F12").Value = "Must order at least 10 at $0.05 each"

# --- Populate the new "Description" text in column E for every part row ---
$ws.Range("E2").Value  = "Microcontroller for rover"
$ws.Range("E3").Value  = "Accelerometer to sense takeoff and landing"
$ws.Range("E4").Value  = "Ultrasonic distance sensor to check for landing"
$ws.Range("E5").Value  = "Reflectance sensor placed on tread 1 to sense distance traveled"
$ws.Range("E6").Value  = "Reflectance sensor placed on tread 2 to sense distance traveled"
$ws.Range("E7").Value  = "Magnetometer to sense direction of rover"
$ws.Range("E8").Value  = "Solenoid to release first marker "
$ws.Range("E9").Value  = "Solenoid to release second marker"
$ws.Range("E10").Value = "Resistor to limit current on transistor for solenoid 1 activation"
$ws.Range("E11").Value = "Resistor to limit current on transistor for solenoid 2 activation"
$ws.Range("E12").Value = "Diode to prevent flyback voltage from solenoid 1"
$ws.Range("E13").Value = "Diode to prevent flyback voltage from solenoid 2"
$ws.Range("E14").Value = "Transistor to activate solenoid 1 from arduino digital output"
$ws.Range("E15").Value = "Transistor to activate solenoid 2 from arduino digital output"

# --- Column widths (closest achievable values given the host's pixel quantization) ---
$ws.Range("A1").ColumnWidth = 10.451822916666666
$ws.Range("B1").ColumnWidth = 19.307291666666668
$ws.Range("C1").ColumnWidth = 11.877604166666666
$ws.Range("D1").ColumnWidth = 87.30729166666667
$ws.Range("E1").ColumnWidth = 61.736979166666664
$ws.Range("F1").ColumnWidth = 37.022135416666664

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("D26").Select() | Out-Null
